$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.135.81"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.783.19"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +1.18%  "
$ws.Range("D5").Value = "'335.31"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").Value = "'0.3779"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.3423"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "'48.29"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "'1.190"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "'0.07423"
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("D12").Value = "'1.007"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "'21.69"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").Value = "'6.420"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Value = "1.784.74"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "'7.039"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "'0.00001089"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "'0.06676"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").Value = "'84.23"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "'6.516"
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").Value = "'17.27"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "27.157.71"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").Value = "'12.37"
$ws.Range("E24").Value = "  -4.99%  "
$ws.Range("D25").Value = "'2.417"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").Value = "'1.508"
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").Value = "'2.527"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("D28").Value = "'21.21"
$ws.Range("E28").Value = "  +6.08%  "
$ws.Range("D29").Value = "'152.77"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "1.986.59"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").Value = "'133.13"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").Value = "'4.058"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "'6.014"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").Value = "'0.08609"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "'12.98"
$ws.Range("E35").Value = "  -1.74%  "
$ws.Range("D36").Value = "'1.651"
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("D37").Value = "'5.417"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("D38").Value = "'0.6807"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").Value = "'0.06345"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").Value = "'8.765"
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("D41").Value = "'0.2183"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").Value = "'0.02327"
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("D43").Value = "'1.256"
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("D44").Value = "'14.44"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").Value = "'0.6372"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "'3.840"
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("D48").Value = "'2.115"
$ws.Range("E48").Value = "  -0.35%  "
$ws.Range("D49").Value = "'128.68"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").Value = "'0.07161"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").Value = "'78.94"
$ws.Range("E51").Value = "  -0.81%  "
